# #114 - Introduce Remarks column on Sales Per Customer report.
# Minor UI enhancements on Reports details page.
#
# The report sheet currently ends at column F ("Bad Order"). We add a new
# "Remarks" column in G, give its header the same look as the other
# header cells (bold, centered - matching D1/E1/F1), and make the column's
# default alignment centered as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default the new column to centered alignment (adds a plain, centered
# cell style used as column G's baseline format).
$ws.Columns("G:G").HorizontalAlignment = -4108

# Write the new header text and copy the existing header formatting
# (bold + centered, same style as the "Bad Order" header in F1) onto it.
$ws.Range("G1").Value = "Remarks"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave column G selected, matching the authored edit.
$ws.Columns("G:G").Select() | Out-Null
